$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "average_doctor" / "average_doctor_old" header labels (BP1/BQ1),
# reflecting that a new "average_doctor" (Harvard case classification) column
# was inserted and the former data shifted into the "_old" slot.
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Updated numeric stats: new "average_doctor" figures in BP, with the
# previous BP values carried over into BQ (average_doctor_old), plus the
# refreshed *_old figures for Ada, Avey, Babylon, Buoy, K health, WebMD,
# doctor_MA, doctor_NJ and doctor_TH across all stat rows.
$changes = @(
    @("E4", 0.485),
    @("F4", 0.054),
    @("G4", 0.232),
    @("N4", 0.474),
    @("O4", 0.058),
    @("P4", 0.242),
    @("Q4", 0.034),
    @("R4", 0.027),
    @("S4", 0.165),
    @("W4", 0.353),
    @("X4", 0.102),
    @("Y4", 0.319),
    @("AI4", 0.423),
    @("AJ4", 0.093),
    @("AK4", 0.306),
    @("AU4", 0.233),
    @("AV4", 0.027),
    @("AW4", 0.165),
    @("BA4", 2.06),
    @("BB4", 0.138),
    @("BC4", 0.371),
    @("BG4", 0.735),
    @("BH4", 0.138),
    @("BI4", 0.371),
    @("BM4", 0.744),
    @("BN4", 0.064),
    @("BO4", 0.253),
    @("BP4", 0.687),
    @("BQ4", 0.757),
    @("E5", 0.609),
    @("F5", 0.062),
    @("G5", 0.249),
    @("N5", 0.729),
    @("O5", 0.081),
    @("P5", 0.284),
    @("Q5", 0.014),
    @("R5", 0.004),
    @("S5", 0.061),
    @("W5", 0.346),
    @("X5", 0.104),
    @("Y5", 0.323),
    @("AI5", 0.428),
    @("AJ5", 0.089),
    @("AK5", 0.299),
    @("AU5", 0.437),
    @("AV5", 0.081),
    @("AW5", 0.285),
    @("BA5", 1.334),
    @("BB5", 0.074),
    @("BG5", 0.387),
    @("BH5", 0.05),
    @("BI5", 0.224),
    @("BM5", 0.538),
    @("BN5", 0.047),
    @("BO5", 0.218),
    @("BP5", 0.445),
    @("BQ5", 0.46),
    @("E6", 0.54),
    @("N6", 0.574),
    @("Q6", 0.02),
    @("W6", 0.349),
    @("AI6", 0.425),
    @("AU6", 0.304),
    @("BA6", 1.611),
    @("BG6", 0.507),
    @("BM6", 0.624),
    @("BP6", 0.537),
    @("BQ6", 0.57),
    @("E7", 0.579),
    @("N7", 0.658),
    @("Q7", 0.016),
    @("W7", 0.347),
    @("AI7", 0.427),
    @("AU7", 0.372),
    @("BA7", 1.432),
    @("BG7", 0.427),
    @("BM7", 0.57),
    @("BP7", 0.477),
    @("BQ7", 0.498),
    @("E8", 0.705),
    @("F8", 0.073),
    @("G8", 0.27),
    @("N8", 0.807),
    @("O8", 0.063),
    @("P8", 0.251),
    @("Q8", 0.02),
    @("S8", 0.108),
    @("W8", 0.39),
    @("X8", 0.119),
    @("Y8", 0.344),
    @("AI8", 0.499),
    @("AJ8", 0.132),
    @("AK8", 0.364),
    @("AU8", 0.378),
    @("AV8", 0.088),
    @("AW8", 0.297),
    @("BA8", 1.791),
    @("BB8", 0.106),
    @("BC8", 0.325),
    @("BG8", 0.568),
    @("BH8", 0.109),
    @("BI8", 0.33),
    @("BM8", 0.681),
    @("BN8", 0.06),
    @("BO8", 0.246),
    @("BP8", 0.597),
    @("BQ8", 0.626),
    @("E9", 0.667),
    @("F9", 0.222),
    @("G9", 0.471),
    @("N9", 0.744),
    @("O9", 0.191),
    @("P9", 0.437),
    @("W9", 0.256),
    @("X9", 0.191),
    @("Y9", 0.437),
    @("AI9", 0.436),
    @("AJ9", 0.246),
    @("AK9", 0.496),
    @("BA9", 1.743),
    @("BB9", 0.25),
    @("BC9", 0.5),
    @("BG9", 0.615),
    @("BH9", 0.237),
    @("BI9", 0.487),
    @("BM9", 0.641),
    @("BN9", 0.23),
    @("BO9", 0.48),
    @("BP9", 0.581),
    @("BQ9", 0.614),
    @("E10", 0.821),
    @("F10", 0.147),
    @("G10", 0.384),
    @("N10", 0.923),
    @("O10", 0.071),
    @("P10", 0.266),
    @("W10", 0.487),
    @("X10", 0.25),
    @("Y10", 0.5),
    @("AI10", 0.538),
    @("AU10", 0.385),
    @("AV10", 0.237),
    @("AW10", 0.487),
    @("BA10", 2.206),
    @("BB10", 0.202),
    @("BC10", 0.45),
    @("BG10", 0.667),
    @("BH10", 0.222),
    @("BI10", 0.471),
    @("BM10", 0.821),
    @("BN10", 0.147),
    @("BO10", 0.384),
    @("BP10", 0.735),
    @("BQ10", 0.765),
    @("E11", 0.846),
    @("F11", 0.13),
    @("G11", 0.361),
    @("N11", 0.923),
    @("O11", 0.071),
    @("P11", 0.266),
    @("W11", 0.487),
    @("X11", 0.25),
    @("Y11", 0.5),
    @("AI11", 0.615),
    @("AJ11", 0.237),
    @("AK11", 0.487),
    @("AU11", 0.513),
    @("AV11", 0.25),
    @("AW11", 0.5),
    @("BA11", 2.206),
    @("BB11", 0.202),
    @("BC11", 0.45),
    @("BG11", 0.667),
    @("BH11", 0.222),
    @("BI11", 0.471),
    @("BM11", 0.821),
    @("BN11", 0.147),
    @("BO11", 0.384),
    @("BP11", 0.735),
    @("BQ11", 0.765),
    @("E12", 1.333),
    @("F12", 0.525),
    @("G12", 0.725),
    @("N12", 1.25),
    @("O12", 0.299),
    @("P12", 0.546),
    @("W12", 1.579),
    @("X12", 0.454),
    @("Y12", 0.674),
    @("AI12", 1.625),
    @("AJ12", 1.484),
    @("AK12", 1.218),
    @("AU12", 2.818),
    @("AV12", 3.512),
    @("AW12", 1.874),
    @("BA12", 3.756),
    @("BB12", 0.459),
    @("BC12", 0.678),
    @("BG12", 1.077),
    @("BH12", 0.071),
    @("BI12", 0.266),
    @("BM12", 1.25),
    @("BN12", 0.25),
    @("BO12", 0.5),
    @("BP12", 1.252),
    @("BQ12", 1.223),
    @("E13", 1.41),
    @("F13", 0.315),
    @("G13", 0.561),
    @("N13", 1.748),
    @("O13", 0.469),
    @("P13", 0.685),
    @("W13", 1.039),
    @("X13", 0.186),
    @("Y13", 0.431),
    @("AI13", 1.174),
    @("AJ13", 0.313),
    @("AK13", 0.56),
    @("AU13", 2.069),
    @("AV13", 0.345),
    @("AW13", 0.587),
    @("BA13", 2.191),
    @("BB13", 0.281),
    @("BC13", 0.531),
    @("BG13", 0.54),
    @("BH13", 0.053),
    @("BI13", 0.229),
    @("BM13", 0.811),
    @("BN13", 0.165),
    @("BO13", 0.406),
    @("BP13", 0.73),
    @("BQ13", 0.671)
)

foreach ($pair in $changes) {
    $addr = $pair[0]
    $val = $pair[1]
    $ws.Range($addr).Value = $val
}
